# [SONALI]: Adding final code with report
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CoverFoxData")

# Update the phone-number value in C2 (was "9860539978", now "9999999999").
# Leading apostrophe forces Excel to store it as text (quote-prefixed),
# matching the existing cell style/formatting instead of converting it to a number.
$ws.Range("C2").Value = "'9999999999"

# Move the active selection to C6
$ws.Range("C6").Select()
